$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.000002731583450289686
$ws.Range("E2").Value = 0.000002731583450289686

# Row 3
$ws.Range("D3").Value = 0.9965475465685959
$ws.Range("E3").Value = 0.9965475465685959

# Row 4
$ws.Range("D4").Value = 0.0000000000000000000001536674452261432
$ws.Range("E4").Value = 0.0000000000000000000001536674452261432

# Row 5
$ws.Range("D5").Value = 0.999702717093282
$ws.Range("E5").Value = 0.999702717093282

# Row 6
$ws.Range("D6").Value = 0.9999803536737446
$ws.Range("E6").Value = 0.9999803536737446

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 4.104519844055176
